# PlayerPerformance_4819.xlsx update:
#  - insert a new "Player Info" sheet before "ODI Batting"
#  - rename MATCH_CARD_LINK -> MATCH_CODE on the batting/bowling sheets and
#    replace the full scorecard URLs with the bare numeric match code

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "Player Info" sheet, inserted before "ODI Batting"
# ---------------------------------------------------------------------------
$infoSheet = $wb.Worksheets.Add($wb.Worksheets.Item("ODI Batting"))
$infoSheet.Name = "Player Info"

$infoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($col = 1; $col -le $infoHeaders.Length; $col++) {
    $cell = $infoSheet.Cells.Item(1, $col)
    $cell.Value = $infoHeaders[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$infoSheet.Cells.Item(2, 1).NumberFormat = "@"
$infoSheet.Cells.Item(2, 1).Value = "4819"
$infoSheet.Cells.Item(2, 2).Value = "Oliver Peter Stone"
$infoSheet.Cells.Item(2, 3).Value = "Right Handed"
$infoSheet.Cells.Item(2, 4).Value = "Right Arm Fast"

# ---------------------------------------------------------------------------
# 2. "ODI Batting" - column D: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Cells.Item(1, 4).Value = "MATCH_CODE"

$battingCodes = @("4209", "4210", "4211", "4212", "4660", "4666", "4698", "4699")
for ($i = 0; $i -lt $battingCodes.Length; $i++) {
    $row = $i + 2
    $cell = $battingSheet.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $battingCodes[$i]
}

# ---------------------------------------------------------------------------
# 3. "ODI Bowling" - column B: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Cells.Item(1, 2).Value = "MATCH_CODE"

$bowlingCodes = @("4210", "4211", "4212", "4660", "4666", "4698", "4699")
for ($i = 0; $i -lt $bowlingCodes.Length; $i++) {
    $row = $i + 2
    $cell = $bowlingSheet.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $bowlingCodes[$i]
}
